$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = -12.501
$ws.Range("C6").Value = -12.318
$ws.Range("C7").Value = -12.909
$ws.Range("E7").Value = 16.385
$ws.Range("E12").Value = 17.498
$ws.Range("E15").Value = 16.172
$ws.Range("C16").Value = -13.319
$ws.Range("C20").Value = -12.205
$ws.Range("E20").Value = 16.337
$ws.Range("E21").Value = 16.534
$ws.Range("E22").Value = 16.453
$ws.Range("E23").Value = 16.434
$ws.Range("C28").Value = -12.848
$ws.Range("C29").Value = -12.129
$ws.Range("E29").Value = 17.002
$ws.Range("C32").Value = -13.057
$ws.Range("E34").Value = 16.849
$ws.Range("C40").Value = -12.151
$ws.Range("E42").Value = 16.539
$ws.Range("E43").Value = 17.051
$ws.Range("E44").Value = 16.465
$ws.Range("E45").Value = 16.544
$ws.Range("C46").Value = -13.654
$ws.Range("E46").Value = 16.819
$ws.Range("E50").Value = 16.433
$ws.Range("C51").Value = -11.257
$ws.Range("E51").Value = 17.157
$ws.Range("C52").Value = -11.355
$ws.Range("C57").Value = -13.94
$ws.Range("C59").Value = -12.918
$ws.Range("C62").Value = -13.555
$ws.Range("C66").Value = -11.579
$ws.Range("E66").Value = 17.157
$ws.Range("E67").Value = 17.314
$ws.Range("C73").Value = -12.332
$ws.Range("C74").Value = -11.866
$ws.Range("E79").Value = 16.862
$ws.Range("E84").Value = 16.406
$ws.Range("C92").Value = -11.066
$ws.Range("E92").Value = 17.769
$ws.Range("E97").Value = 16.77
$ws.Range("C100").Value = -12.726
